$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue $ws "D2" "237.39"
Set-TextValue $ws "D3" "22.46"
Set-TextValue $ws "D4" "5.406"
Set-TextValue $ws "D5" "0.05645"
Set-TextValue $ws "D7" "6.479"
Set-TextValue $ws "D8" "1.079"
Set-TextValue $ws "B10" "WazirX"
Set-TextValue $ws "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws "D10" "0.1400"
Set-TextValue $ws "E10" "9WazirXWRX"
Set-TextValue $ws "B11" "MandalaExchangeToken"
Set-TextValue $ws "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws "D11" "0.07321"
Set-TextValue $ws "E11" "10MandalaExchangeTokenMDX"
Set-TextValue $ws "B12" "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws "D12" "0.03204"
Set-TextValue $ws "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue $ws "B13" "BitrueCoin"
Set-TextValue $ws "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws "D13" "0.02936"
Set-TextValue $ws "E13" "12BitrueCoinBTR"
Set-TextValue $ws "B14" "BitMartToken"
Set-TextValue $ws "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws "D14" "0.09259"
Set-TextValue $ws "E14" "13BitMartTokenBMX"
Set-TextValue $ws "B15" "BitForexToken"
Set-TextValue $ws "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws "D15" "0.001672"
Set-TextValue $ws "E15" "14BitForexTokenBF"
Set-TextValue $ws "B16" "MCDex"
Set-TextValue $ws "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws "D16" "3.252"
Set-TextValue $ws "E16" "15MCDexMCB"
Set-TextValue $ws "B17" "CoinExToken"
Set-TextValue $ws "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws "D17" "0.04755"
Set-TextValue $ws "E17" "16CoinExTokenCET"
Set-TextValue $ws "B18" "One"
Set-TextValue $ws "C18" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws "D18" "0.0005741"
Set-TextValue $ws "E18" "17OneONE"
Set-TextValue $ws "D19" "0.006231"
Set-TextValue $ws "D20" "0.005108"
Set-TextValue $ws "D21" "0.001051"
Set-TextValue $ws "D22" "0.0001500"
Set-TextValue $ws "D23" "3.860"
Set-TextValue $ws "E23" "22LEOLEOBestin24h"
Set-TextValue $ws "D26" "0.1056"
Set-TextValue $ws "D27" "0.0004991"
Set-TextValue $ws "D40" "0.04078"
Set-TextValue $ws "D41" "0.006981"
Set-TextValue $ws "B42" "BKEXToken"
Set-TextValue $ws "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D42" "0.1041"
Set-TextValue $ws "E42" "41BKEXTokenBKK"
Set-TextValue $ws "B43" "CEJI"
Set-TextValue $ws "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.003209"
Set-TextValue $ws "E43" "42CEJICEJI"
Set-TextValue $ws "D44" "0.009901"
Set-TextValue $ws "D45" "0.00005414"
Set-TextValue $ws "D47" "0.6754"
Set-TextValue $ws "D48" "0.03917"
